{"js": "// Update the date line and all 25 two-digit-by-two-digit multiplication\n// problems in the table to the new values from the next day's worksheet.\nconst replacements = [\n  [\"2024-09-26 Thursday\", \"2024-09-27 Friday\"],\n  [\"75\u00d731=\", \"73\u00d715=\"],\n  [\"79\u00d779=\", \"23\u00d711=\"],\n  [\"26\u00d743=\", \"28\u00d739=\"],\n  [\"68\u00d736=\", \"60\u00d736=\"],\n  [\"51\u00d739=\", \"67\u00d739=\"],\n  [\"35\u00d773=\", \"36\u00d783=\"],\n  [\"64\u00d757=\", \"93\u00d791=\"],\n  [\"24\u00d785=\", \"21\u00d712=\"],\n  [\"53\u00d739=\", \"37\u00d758=\"],\n  [\"93\u00d747=\", \"45\u00d733=\"],\n  [\"88\u00d721=\", \"79\u00d797=\"],\n  [\"46\u00d749=\", \"24\u00d744=\"],\n  [\"40\u00d782=\", \"34\u00d776=\"],\n  [\"51\u00d735=\", \"63\u00d771=\"],\n  [\"45\u00d759=\", \"92\u00d738=\"],\n  [\"72\u00d776=\", \"16\u00d717=\"],\n  [\"61\u00d775=\", \"89\u00d786=\"],\n  [\"87\u00d769=\", \"25\u00d724=\"],\n  [\"99\u00d759=\", \"95\u00d723=\"],\n  [\"91\u00d725=\", \"14\u00d722=\"],\n  [\"37\u00d773=\", \"25\u00d761=\"],\n  [\"21\u00d718=\", \"51\u00d732=\"],\n  [\"47\u00d789=\", \"21\u00d725=\"],\n  [\"55\u00d782=\", \"26\u00d714=\"],\n  [\"30\u00d746=\", \"51\u00d780=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and all 25 two-digit-by-two-digit multiplication\n# problems in the table to the new values from the next day's worksheet.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-09-26 Thursday\", \"2024-09-27 Friday\"),\n    @(\"75\u00d731=\", \"73\u00d715=\"),\n    @(\"79\u00d779=\", \"23\u00d711=\"),\n    @(\"26\u00d743=\", \"28\u00d739=\"),\n    @(\"68\u00d736=\", \"60\u00d736=\"),\n    @(\"51\u00d739=\", \"67\u00d739=\"),\n    @(\"35\u00d773=\", \"36\u00d783=\"),\n    @(\"64\u00d757=\", \"93\u00d791=\"),\n    @(\"24\u00d785=\", \"21\u00d712=\"),\n    @(\"53\u00d739=\", \"37\u00d758=\"),\n    @(\"93\u00d747=\", \"45\u00d733=\"),\n    @(\"88\u00d721=\", \"79\u00d797=\"),\n    @(\"46\u00d749=\", \"24\u00d744=\"),\n    @(\"40\u00d782=\", \"34\u00d776=\"),\n    @(\"51\u00d735=\", \"63\u00d771=\"),\n    @(\"45\u00d759=\", \"92\u00d738=\"),\n    @(\"72\u00d776=\", \"16\u00d717=\"),\n    @(\"61\u00d775=\", \"89\u00d786=\"),\n    @(\"87\u00d769=\", \"25\u00d724=\"),\n    @(\"99\u00d759=\", \"95\u00d723=\"),\n    @(\"91\u00d725=\", \"14\u00d722=\"),\n    @(\"37\u00d773=\", \"25\u00d761=\"),\n    @(\"21\u00d718=\", \"51\u00d732=\"),\n    @(\"47\u00d789=\", \"21\u00d725=\"),\n    @(\"55\u00d782=\", \"26\u00d714=\"),\n    @(\"30\u00d746=\", \"51\u00d780=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
